# Update the PGT doc-tracking sheet: 5 new rows were recorded for this export,
# shifting/reflowing the "Desconhecido"/"ANDER RODOLFO HENRIQUE"/"ROSELI NUNES"
# block (rows 163-201) down into rows 163-206.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the 5 additional rows (new dimension becomes A1:F206).
for ($i = 0; $i -lt 5; $i++) {
    $ws.Rows(163).Insert()
}

# Final values for row 163 through 206 (Tipo de documento PGT, Assentamento,
# Municipio, Nome T1, Autenticador, Objetivo).
$rowData = @{
    163 = @("Relatório de conformidades para regularização", "Desconhecido", "Desconhecido", "Valderi", "D059819FE574FBA27BD1084144C2533C", "Regularização")
    164 = @("Solicitação de documentação complementar", "ACOPAM", "BOA VENTURA DE SAO ROQUE", "Valderi", "A963B2E2DC5B0F6BD29DC8A31B09D4EF", "")
    165 = @("Relatório de conformidades para regularização", "Desconhecido", "Desconhecido", "Adelar", "98518C564982BC7968097BDEDA2A9BDB", "Regularização")
    166 = @("Solicitação de documentação complementar", "Desconhecido", "Desconhecido", "Adelar", "3FF2D22819B3330B814BC9017012393F", "")
    167 = @("Análise para regularização", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Sandra", "112D6F0C94649FCED3C07DC35802FB79", "")
    168 = @("Relatório de conformidades para regularização", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Sandra", "3BC79CA16E0DD94D83CCEF5BE4EFD256", "Regularização")
    169 = @("Análise para regularização", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Simone", "183830569A2CEF2EEF133D2A68AB051A", "")
    170 = @("Relatório de conformidades para regularização", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Simone", "866FFD0DE100FC917059EF68503DA039", "Regularização")
    171 = @("Relatório de conformidades para regularização", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Simone", "B68A9D35C6E7C78C65B895D9B4E8DCBA", "Regularização")
    172 = @("Solicitação de documentação complementar", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Simone", "DB62A611FC09102671B3EA905E42C24C", "")
    173 = @("Relatório de conformidades para regularização", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Andiara", "1E28611DB40D7BCF9AF123CC4CC60BFE", "Regularização")
    174 = @("Solicitação de documentação complementar", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Andiara", "069ACF289074156352A4781DFCD0D264", "")
    175 = @("Análise para regularização", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Cristiane", "180B234A9A4D11A43083BF6E6BFF3617", "")
    176 = @("Relatório de conformidades para regularização", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Cristiane", "5AB0DA0076CB4A7BC654C14601170805", "Regularização")
    177 = @("Relatório de conformidades para regularização", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Cristiane", "B44DAF5094FF26D7D29600F23051E9C7", "Regularização")
    178 = @("Solicitação de documentação complementar", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Cristiane", "0B020F956312AEDBA007ECBA724ED40B", "")
    179 = @("Relatório de conformidades para regularização", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Claudineia", "58258FAC18CC8DE705C5107FC0A0B759", "Regularização")
    180 = @("Solicitação de documentação complementar", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Claudineia", "AF93B4687B22186B9F40FA9CAF59355D", "")
    181 = @("Análise para regularização", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Milene", "5D79082732416772252230E1D5490EFA", "")
    182 = @("Relatório de conformidades para regularização", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Milene", "55DD509716CC263B93E3CBE427D08A02", "Regularização")
    183 = @("Análise para regularização", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Luciane", "C54F6E8F93CF6FDCB460C4EC27CB5D0B", "")
    184 = @("Relatório de conformidades para regularização", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Luciane", "56FEDF1CFCC0CE78021FA9F765943DA5", "Regularização")
    185 = @("Relatório de conformidades para regularização", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Luciane", "6B45729139C2A27BC5092915C751010D", "Regularização")
    186 = @("Solicitação de documentação complementar", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Luciane", "2A8C3895F9792659241327C30A81AC6D", "")
    187 = @("Análise para regularização", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Gabriel", "962E98FB0F61D261DB0D48C557F70ED2", "")
    188 = @("Relatório de conformidades para regularização", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Gabriel", "42A275CFB3A429A9CA6E7AB9FB7521B6", "Regularização")
    189 = @("Relatório de conformidades para regularização", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Gabriel", "9D159FDA14EB647768A2EF8FCA866AF5", "Regularização")
    190 = @("Solicitação de documentação complementar", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Gabriel", "912B181376702376839E92A34234AB6F", "")
    191 = @("Relatório de conformidades para regularização", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Marlei", "1A97F4E202A341184318AE16C57D676A", "Regularização")
    192 = @("Solicitação de documentação complementar", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Marlei", "7ACF437DA3089897B5DEF97EE5ADB295", "")
    193 = @("Relatório de conformidades para regularização", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Laudiceia", "14C5879842FBF191CCE378771DE20033", "Regularização")
    194 = @("Solicitação de documentação complementar", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Vilmar", "F77DBF14F49FCBA5D570864B9B5782D0", "")
    195 = @("Relatório de conformidades para regularização", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Juliana", "916E2681C0054750A9E687574A0B6489", "Regularização")
    196 = @("Solicitação de documentação complementar", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Juliana", "F17621921D038A38C3EB0D53E3733F98", "")
    197 = @("Relatório de conformidades para regularização", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Rosali", "1981ACD0B95BCA268D57836FFB95AE3A", "Regularização")
    198 = @("Solicitação de documentação complementar", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Rosali", "C154AC8381DB0FCD102A1559720175A7", "")
    199 = @("Relatório de conformidades para regularização", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Peterson", "BC998959680EB2D70139D8E84542E2F2", "Regularização")
    200 = @("Solicitação de documentação complementar", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Peterson", "EFDC3F40D59FE406CC2DEF20D8044981", "")
    201 = @("Relatório de conformidades para regularização", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Adrian", "7A7091058DBBDBDBFFEC38A46A34A73F", "Regularização")
    202 = @("Solicitação de documentação complementar", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Adrian", "7DE34ABB18DFC4184EA4ECFE7822F34E", "")
    203 = @("Relatório de conformidades para regularização", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Debora", "874465517EEF60A1E4B3C3E4DBF99215", "Regularização")
    204 = @("Solicitação de documentação complementar", "ANDER RODOLFO HENRIQUE", "DIAMANTE DO OESTE", "Debora", "25B5DA7A79DCC134ADF32C921CBAB522", "")
    205 = @("Análise para regularização", "ROSELI NUNES", "SAO JERONIMO DA SERRA", "Guilherme", "4C91EE4CEE2814F5667059FBECB733FC", "")
    206 = @("Relatório de conformidades para regularização", "ROSELI NUNES", "SAO JERONIMO DA SERRA", "Guilherme", "B53F099B4CC28D9A9292715B93005CED", "Regularização")
}

foreach ($rowNum in $rowData.Keys) {
    $values = $rowData[$rowNum]
    for ($col = 1; $col -le $values.Length; $col++) {
        $ws.Cells.Item([int]$rowNum, $col).Value = $values[$col - 1]
    }
}
